$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 header: "PLU"/"SERIAL" columns are repurposed to "Vendedor"/"Cedula
# Cliente". Columns C/D (MSIDN/MSI) stay as-is.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value2 = "Vendedor"
$ws.Range("B10").Value2 = "Cedula Cliente"

# Give A10:B10 the same bold/centered look as before, but re-touch the fill
# so the cell format gets its own explicit (empty) fill applied - matching
# the restyled header cells in the updated sheet.
$hdr = $ws.Range("A10:B10")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.NumberFormat = "@"
$hdr.Interior.ColorIndex = -4142

# ---------------------------------------------------------------------------
# Data rows 11-13: new Vendedor/Cedula Cliente values paired with the
# existing MSIDN/MSI sample rows (row 12 keeps the second MSIDN/MSI pair,
# and row 13 repeats the first pair).
# ---------------------------------------------------------------------------
$ws.Range("A11").Value2 = "10960370"
$ws.Range("B11").Value2 = "667299000"
$ws.Range("C11").Value2 = "3016875982"
$ws.Range("D11").Value2 = "732111198172291"

$ws.Range("A12").Value2 = "10960370"
$ws.Range("B12").Value2 = "667299000"
$ws.Range("C12").Value2 = "3016875893"
$ws.Range("D12").Value2 = "732111198172290"

$ws.Range("A13").Value2 = "10960370"
$ws.Range("B13").Value2 = "667299000"
$ws.Range("C13").Value2 = "3016875982"
$ws.Range("D13").Value2 = "732111198172291"

# Vendedor / Cedula Cliente columns (A:B) on the data rows get their own
# text style with an explicit (empty) fill, distinct from the plain text
# style still used by the MSIDN/MSI columns (C:D).
$data = $ws.Range("A11:B13")
$data.NumberFormat = "@"
$data.Interior.ColorIndex = -4142

$msi = $ws.Range("C11:D13")
$msi.NumberFormat = "@"

# ---------------------------------------------------------------------------
# View state: selection moves to D14 and the previous frozen/scrolled
# top-left cell is cleared.
# ---------------------------------------------------------------------------
$ws.Range("D14").Select()
